$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.097.64"
$ws.Range("E2").Value = "  -3.13%  "
$ws.Range("D3").Value = "'3.803.92"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'594.98"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").Value = "'172.90"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("D7").Value = "'3.802.59"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").Value = "'6.24"
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("D13").Value = "'38.03"
$ws.Range("E13").Value = "  -5.50%  "
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "'4.438.60"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "'3.800.86"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "'68.225.11"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").Value = "'7.15"
$ws.Range("E19").Value = "  -5.41%  "
$ws.Range("D20").Value = "'16.08"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "'488.61"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").Value = "'9.27"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "'0.734"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "'84.62"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -8.68%  "
$ws.Range("E26").Value = "  +3.14%  "
$ws.Range("E27").Value = "  -6.68%  "
$ws.Range("E28").Value = "  -9.55%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").Value = "'32.78"
$ws.Range("E32").Value = "  +6.88%  "
$ws.Range("D33").Value = "'7.73"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "'5.78"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("D39").Value = "'0.326"
$ws.Range("E39").Value = "  -7.52%  "
$ws.Range("D40").Value = "'455.55"
$ws.Range("E40").Value = "  +5.00%  "
$ws.Range("D41").Value = "'48.89"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").Value = "'2.89"
$ws.Range("E43").Value = "  -6.79%  "
$ws.Range("D44").Value = "'8.27"
$ws.Range("E44").Value = "  -4.34%  "
$ws.Range("D45").Value = "'41.49"
$ws.Range("E45").Value = "  -8.87%  "
$ws.Range("D46").Value = "'2.828.17"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").Value = "'139.94"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "'0.0351"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").Value = "'26.33"
$ws.Range("E50").Value = "  -4.37%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  -6.90%  "
